$wb = $excel.ActiveWorkbook

# --- LAPSE (0->4) sheet: insert a new row above row 5 (new CALENDARYEAR 2021) ---
$wsLapse = $wb.Worksheets.Item("LAPSE (0->4)")
[void]$wsLapse.Rows("5:5").Insert(-4121, 0)

# New row 5: year 2021 with the same 5% lapse rate / formatting as the old row 5
[void]$wsLapse.Range("A5").ClearFormats()
$wsLapse.Range("A5").Value = 2021
$wsLapse.Range("B5").Value = 0.05
$wsLapse.Range("B5").NumberFormat = "0%"

# Move the selection to A6 and make this sheet the active tab
[void]$wsLapse.Range("A6").Select()
[void]$wsLapse.Activate()

# --- DIS2(0->2) sheet: scroll the view back to the top-left (drop topLeftCell) ---
$wsDis2 = $wb.Worksheets.Item("DIS2(0->2)")
[void]$wsDis2.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1

# restore LAPSE (0->4) as the active / selected sheet
[void]$wsLapse.Activate()
